$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$rng = $ws.Range("A1:U80")
$tbl = $ws.ListObjects.Add(1, $rng, [System.Reflection.Missing]::Value, 1)
try {
    $tbl.TableStyle = "TableStyleNone"
    Write-Output ("ok: " + $tbl.TableStyle.Name)
} catch {
    Write-Output "err1: $_"
}
